# Automatische test-sync: 2025-08-30 18:39:50
# Adds a new "Afspraak demo" mail-log entry (row 4) to the Logs sheet,
# updates the Dashboard summary with a new "Planning / Afspraak" category
# row (row 3), extends the matching conditional-formatting ranges, and
# widens the chart series ranges so the new Dashboard row is plotted.

$wb = $excel.ActiveWorkbook

# --- Logs sheet: append row 4 -------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Cells.Item(4, 1).Value  = "Afspraak demo"
$logs.Cells.Item(4, 2).Value  = "mailmind.test@zohomail.eu"
$logs.Cells.Item(4, 4).Value  = "Planning / Afspraak"
$logs.Cells.Item(4, 6).Value  = "2025-08-30 18:39:20"
$logs.Cells.Item(4, 7).Value  = "Nee"
$logs.Cells.Item(4, 8).Value  = "Ja"
$logs.Cells.Item(4, 9).Value  = "Nee"
$logs.Cells.Item(4, 10).Value = "Nee"

# Extend the conditional-formatting ranges (D/G/H/I/J) from row 2:3 to 2:4,
# keeping their existing rules intact.
$oldRanges = @("D2:D3", "G2:G3", "H2:H3", "I2:I3", "J2:J3")
$newRanges = @("D2:D4", "G2:G4", "H2:H4", "I2:I4", "J2:J4")
for ($i = 0; $i -lt $oldRanges.Length; $i++) {
    $fcs = $logs.Range($oldRanges[$i]).FormatConditions
    for ($j = 1; $j -le $fcs.Count; $j++) {
        $fcs.Item($j).ModifyAppliesToRange($logs.Range($newRanges[$i]))
    }
}

# --- Dashboard sheet: append row 3 ---------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Cells.Item(3, 1).Value = "Planning / Afspraak"
$dash.Cells.Item(3, 2).Value = 1

# --- Chart: extend series ranges to include the new Dashboard row -------
# Assign string formulas (not Range objects) so only the cat/val refs move
# and the series name (tx) reference is left untouched.
$chart = $dash.ChartObjects(1).Chart
$series = $chart.SeriesCollection(1)
$series.XValues = "='Dashboard'!`$A`$2:`$A`$3"
$series.Values = "='Dashboard'!`$B`$2:`$B`$3"
